$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 643.3333
$ws.Range("I15").Value = 643.3333
$ws.Range("K15").Value = 1929.9999
$ws.Range("M15").Value = -1760.9999
$ws.Range("H42").Value = 2114
$ws.Range("I42").Value = 125.28571
$ws.Range("K42").Value = 375.85713
$ws.Range("M42").Value = -145.85713
$ws.Range("H43").Value = 9499.5
$ws.Range("I43").Value = 9000
$ws.Range("K43").Value = 9000
$ws.Range("M43").Value = -8931
$ws.Range("H80").Value = 1801.25
$ws.Range("I80").Value = 3591.2
$ws.Range("K80").Value = 10773.6
$ws.Range("M80").Value = -9775.599999999999
$ws.Range("H83").Value = 1801.25
$ws.Range("I83").Value = 3591.2
$ws.Range("K83").Value = 32320.8
$ws.Range("M83").Value = -27328.8
$ws.Range("H88").Value = 24617440
$ws.Range("I88").Value = 111119450
$ws.Range("J88").Value = 2991939
$ws.Range("K88").Value = 111119450
$ws.Range("L88").Value = 2991939
$ws.Range("M88").Value = -111119044
$ws.Range("N88").Value = -2992751
$ws.Range("H91").Value = 24617440
$ws.Range("I91").Value = 111119450
$ws.Range("J91").Value = 2991939
$ws.Range("K91").Value = 111119450
$ws.Range("L91").Value = 2991939
$ws.Range("M91").Value = -111118046
$ws.Range("N91").Value = -2994747
$ws.Range("H100").Value = 2363.111
$ws.Range("I100").Value = 2826.8572
$ws.Range("J100").Value = 740
$ws.Range("K100").Value = 2826.8572
$ws.Range("L100").Value = 740
$ws.Range("M100").Value = -2285.8572
$ws.Range("N100").Value = -1822
$ws.Range("H106").Value = 2546.1333
$ws.Range("I106").Value = 2777
$ws.Range("J106").Value = 2199.8333
$ws.Range("K106").Value = 2777
$ws.Range("L106").Value = 2199.8333
$ws.Range("M106").Value = -2146
$ws.Range("N106").Value = -3461.8333
$ws.Range("H132").Value = 3927.4285
$ws.Range("I132").Value = 4415.3335
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 13246.0005
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -10716.0005
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 2656.6956
$ws.Range("J137").Value = 3914.1667
$ws.Range("L137").Value = 11742.5001
$ws.Range("N137").Value = -16842.5001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").ClearContents()  # was 24999
$ws.Range("N18").Value = 0
$ws.Range("H61").Value = 31252812
$ws.Range("I61").Value = 45457436
$ws.Range("J61").Value = 2640.7
$ws.Range("K61").Value = 45457436
$ws.Range("L61").Value = 2640.7
$ws.Range("M61").Value = -45457224
$ws.Range("N61").Value = -3064.7
$ws.Range("H136").Value = 31252812
$ws.Range("I136").Value = 45457436
$ws.Range("J136").Value = 2640.7
$ws.Range("K136").Value = 136372308
$ws.Range("L136").Value = 7922.099999999999
$ws.Range("M136").Value = -136369758
$ws.Range("N136").Value = -13022.1

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()  # was -14875
$ws.Range("H99").Value = 2199.6667
$ws.Range("J99").Value = 3199.6667
$ws.Range("L99").Value = 3199.6667
$ws.Range("N99").Value = -6195.6667
$ws.Range("H134").Value = 15156944
$ws.Range("I134").Value = 16672280
$ws.Range("J134").Value = 3583
$ws.Range("K134").Value = 50016840
$ws.Range("L134").Value = 10749
$ws.Range("M134").Value = -50014305
$ws.Range("N134").Value = -15819

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4970.7085
$ws.Range("I31").Value = 3405.389
$ws.Range("J31").Value = 9666.666999999999
$ws.Range("K31").Value = 3405.389
$ws.Range("L31").Value = 9666.666999999999
$ws.Range("M31").Value = -3110.389
$ws.Range("N31").Value = -10256.667
$ws.Range("H34").Value = 4970.7085
$ws.Range("I34").Value = 3405.389
$ws.Range("J34").Value = 9666.666999999999
$ws.Range("K34").Value = 3405.389
$ws.Range("L34").Value = 9666.666999999999
$ws.Range("M34").Value = -3203.389
$ws.Range("N34").Value = -10070.667
$ws.Range("H122").Value = 1807.6471
$ws.Range("I122").Value = 2040.4166
$ws.Range("K122").Value = 6121.2498
$ws.Range("M122").Value = -3671.2498
$ws.Range("H132").Value = 47622480
$ws.Range("I132").Value = 66669870
$ws.Range("K132").Value = 200009610
$ws.Range("M132").Value = -200007080

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 180235.3
$ws.Range("J11").Value = 64500
$ws.Range("L11").Value = 193500
$ws.Range("N11").Value = -193780
$ws.Range("H97").Value = 423.22223
$ws.Range("I97").Value = 424.75
$ws.Range("J97").Value = 422
$ws.Range("K97").Value = 1274.25
$ws.Range("L97").Value = 1266
$ws.Range("M97").Value = -778.25
$ws.Range("N97").Value = -2258
$ws.Range("H107").Value = 1115.3572
$ws.Range("I107").Value = 467.16666
$ws.Range("J107").Value = 1601.5
$ws.Range("K107").Value = 1401.49998
$ws.Range("L107").Value = 4804.5
$ws.Range("M107").Value = 518.5000199999999
$ws.Range("N107").Value = -8644.5
$ws.Range("H122").Value = 1583.909
$ws.Range("I122").Value = 1431.5
$ws.Range("K122").Value = 12883.5
$ws.Range("M122").Value = -10433.5
$ws.Range("H129").Value = 3002
$ws.Range("I129").Value = 789.6667
$ws.Range("K129").Value = 2369.0001
$ws.Range("M129").Value = 2630.9999
$ws.Range("H133").Value = 10854.8
$ws.Range("I133").Value = 3924.6667
$ws.Range("J133").Value = 21250
$ws.Range("K133").Value = 11774.0001
$ws.Range("L133").Value = 63750
$ws.Range("M133").Value = -6714.000100000001
$ws.Range("N133").Value = -73870
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()  # was 24000
$ws.Range("N139").Value = 0

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2483.3333
$ws.Range("I102").Value = 1981.8182
$ws.Range("J102").Value = 8000
$ws.Range("K102").Value = 1981.8182
$ws.Range("L102").Value = 8000
$ws.Range("M102").Value = -359.8181999999999
$ws.Range("N102").Value = -11244
$ws.Range("H132").Value = 7357004
$ws.Range("I132").Value = 8932648
$ws.Range("K132").Value = 26797944
$ws.Range("M132").Value = -26795414

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 9999.75
$ws.Range("I26").Value = 6666.6665
$ws.Range("K26").Value = 6666.6665
$ws.Range("M26").Value = -6371.6665
$ws.Range("I46").Value = 2427.5
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 2427.5
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -2239.5
$ws.Range("N46").Value = -2376
$ws.Range("H136").Value = 1637.1364
$ws.Range("J136").Value = 2096
$ws.Range("L136").Value = 6288
$ws.Range("N136").Value = -11388

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 41493
$ws.Range("I41").Value = 28999
$ws.Range("J41").Value = 44616.5
$ws.Range("K41").Value = 28999
$ws.Range("L41").Value = 44616.5
$ws.Range("M41").Value = -28609
$ws.Range("N41").Value = -45396.5
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()  # was -4431
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()  # was 44999
$ws.Range("N50").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()  # was 25000
$ws.Range("N51").Value = 0
$ws.Range("H52").Value = 14500
$ws.Range("I52").Value = 14000
$ws.Range("K52").Value = 14000
$ws.Range("M52").Value = -13774
$ws.Range("H58").Value = 42200
$ws.Range("J58").Value = 42200
$ws.Range("L58").Value = 42200
$ws.Range("N58").Value = -42816
$ws.Range("H81").Value = 6200
$ws.Range("I81").Value = 4599.3335
$ws.Range("K81").Value = 9198.666999999999
$ws.Range("M81").Value = -8137.666999999999
$ws.Range("H84").Value = 6200
$ws.Range("I84").Value = 4599.3335
$ws.Range("K84").Value = 45993.335
$ws.Range("M84").Value = -40689.335
$ws.Range("H100").Value = 2043.1666
$ws.Range("I100").Value = 2057.8667
$ws.Range("K100").Value = 4115.7334
$ws.Range("M100").Value = -3574.7334
$ws.Range("H122").Value = 2996
$ws.Range("I122").Value = 2996
$ws.Range("K122").Value = 8988
$ws.Range("M122").Value = -6538
$ws.Range("H126").Value = 2932.9412
$ws.Range("J126").Value = 3450
$ws.Range("L126").Value = 10350
$ws.Range("N126").Value = -15290
$ws.Range("H136").Value = 21741502
$ws.Range("I136").Value = 22729662
$ws.Range("K136").Value = 68188986
$ws.Range("M136").Value = -68186436
